$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new string values in the precise order required so that the
# shared-strings table is appended to in the same order as the source
# workbook: COND_BOTTOM, COND_SURFACE, CONDB, CONDS, TEMPS, TEMPB,
# TEMP_SURFACE, TEMP_BOTTOM.
$ws.Range("A24").Value = "COND_BOTTOM"
$ws.Range("A25").Value = "COND_SURFACE"
$ws.Range("B24").Value = "CONDB"
$ws.Range("B25").Value = "CONDS"
$ws.Range("B26").Value = "TEMPS"
$ws.Range("B27").Value = "TEMPB"
$ws.Range("A26").Value = "TEMP_SURFACE"
$ws.Range("A27").Value = "TEMP_BOTTOM"

# Numeric columns C (Conversion flag) and D (Conversion factor)
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 1

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1

$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 1

$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 1

# Update the active selection to mirror the author's final cursor position
$ws.Range("A27").Select() | Out-Null
